$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the full address locator text: "East" -> "E"
# (leading apostrophe preserves the existing quote-prefix cell style)
$ws.Range("E2").Value = "'1101 E Karsch Blvd, Farmington, MO 63640"

# Update the active selection/cursor position
$ws.Range("E3").Select()
